{"js": "// Update the date line and the 25 division problems in the table to the\n// new values, per the commit diff.\nconst replacements = [\n  [\"2024-02-07 Wednesday\", \"2024-02-08 Thursday\"],\n  [\"536\u00f74=\", \"799\u00f73=\"],\n  [\"510\u00f77=\", \"579\u00f78=\"],\n  [\"237\u00f73=\", \"870\u00f78=\"],\n  [\"112\u00f73=\", \"455\u00f76=\"],\n  [\"630\u00f76=\", \"919\u00f77=\"],\n  [\"968\u00f79=\", \"425\u00f72=\"],\n  [\"301\u00f73=\", \"294\u00f78=\"],\n  [\"263\u00f77=\", \"741\u00f76=\"],\n  [\"965\u00f78=\", \"448\u00f77=\"],\n  [\"481\u00f77=\", \"284\u00f78=\"],\n  [\"690\u00f73=\", \"577\u00f76=\"],\n  [\"881\u00f75=\", \"135\u00f74=\"],\n  [\"471\u00f77=\", \"397\u00f73=\"],\n  [\"379\u00f75=\", \"467\u00f77=\"],\n  [\"299\u00f78=\", \"647\u00f77=\"],\n  [\"698\u00f75=\", \"375\u00f79=\"],\n  [\"819\u00f72=\", \"914\u00f74=\"],\n  [\"378\u00f79=\", \"585\u00f72=\"],\n  [\"179\u00f75=\", \"133\u00f75=\"],\n  [\"900\u00f78=\", \"207\u00f79=\"],\n  [\"400\u00f77=\", \"682\u00f77=\"],\n  [\"180\u00f74=\", \"321\u00f75=\"],\n  [\"652\u00f73=\", \"512\u00f72=\"],\n  [\"851\u00f76=\", \"703\u00f74=\"],\n  [\"484\u00f79=\", \"882\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division problems in the table to the\n# new values, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-07 Wednesday\", \"2024-02-08 Thursday\"),\n    @(\"536\u00f74=\", \"799\u00f73=\"),\n    @(\"510\u00f77=\", \"579\u00f78=\"),\n    @(\"237\u00f73=\", \"870\u00f78=\"),\n    @(\"112\u00f73=\", \"455\u00f76=\"),\n    @(\"630\u00f76=\", \"919\u00f77=\"),\n    @(\"968\u00f79=\", \"425\u00f72=\"),\n    @(\"301\u00f73=\", \"294\u00f78=\"),\n    @(\"263\u00f77=\", \"741\u00f76=\"),\n    @(\"965\u00f78=\", \"448\u00f77=\"),\n    @(\"481\u00f77=\", \"284\u00f78=\"),\n    @(\"690\u00f73=\", \"577\u00f76=\"),\n    @(\"881\u00f75=\", \"135\u00f74=\"),\n    @(\"471\u00f77=\", \"397\u00f73=\"),\n    @(\"379\u00f75=\", \"467\u00f77=\"),\n    @(\"299\u00f78=\", \"647\u00f77=\"),\n    @(\"698\u00f75=\", \"375\u00f79=\"),\n    @(\"819\u00f72=\", \"914\u00f74=\"),\n    @(\"378\u00f79=\", \"585\u00f72=\"),\n    @(\"179\u00f75=\", \"133\u00f75=\"),\n    @(\"900\u00f78=\", \"207\u00f79=\"),\n    @(\"400\u00f77=\", \"682\u00f77=\"),\n    @(\"180\u00f74=\", \"321\u00f75=\"),\n    @(\"652\u00f73=\", \"512\u00f72=\"),\n    @(\"851\u00f76=\", \"703\u00f74=\"),\n    @(\"484\u00f79=\", \"882\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # Use a fresh Find on the whole document content each iteration so\n    # the search range/position is reset and wdReplaceAll finds every\n    # occurrence (there is exactly one of each in this document).\n    $find = $d.Content.Find\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
